$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.921.77'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.75%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.746.50'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.01%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.62%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.37'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.70%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.744.31'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.97%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("E9").Value = '  +2.26%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.167'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.24%  '

$ws.Range("E11").Value = '  +3.77%  '

$ws.Range("E12").Value = '  +0.61%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.32'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.57%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000249'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.08%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.373.43'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.28%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.742.51'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.23%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.946.63'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.82%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.29'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.28%  '

$ws.Range("E19").Value = '  +0.39%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.12'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.38%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.84'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +19.76%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '495.23'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.38%  '

$ws.Range("E23").Value = '  +1.90%  '

$ws.Range("E24").Value = '  +12.44%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.41'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.21%  '

$ws.Range("E26").Value = '  +2.06%  '

$ws.Range("E27").Value = '  +2.40%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.38'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.65%  '

$ws.Range("E29").Value = '  +0.41%  '

$ws.Range("E30").Value = '  +7.44%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.98'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.74%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.96'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.50%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.87'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.03%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.892.49'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.48%  '

$ws.Range("E35").Value = '  +2.20%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.680.89'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.04%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.20%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.02'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.80%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.88'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.61%  '

$ws.Range("E40").Value = '  +1.71%  '

$ws.Range("E41").Value = '  +0.96%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '438.86'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.40%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +8.60%  '

$ws.Range("E44").Value = '  +0.78%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.98'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.87%  '

$ws.Range("E46").Value = '  +2.08%  '

$ws.Range("E47").Value = '  +0.01%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.62'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.46%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.17'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.64%  '

$ws.Range("E50").Value = '  +3.30%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.783.15'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.39%  '
